# EPBDS: added conditionABId to key, updated xls for Class level configuration
#
# This script extends the "Class level configuration" sample table
# (rows 21-24, columns C-J) with two new columns (K: classABeanFactory,
# L: classBBeanFactory) and a new data row (25) that mirrors row 24's
# layout/format but uses the "A"/"C" condition-key pairing (like row 16
# uses "A"/"C" for the Field level table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# D/E/J shrink slightly (and lose their "best fit" flag, same as Excel
# does whenever a width is set explicitly); K and L are brand new columns
# that hold the two new "Bean Factory" headers.
$ws.Columns.Item(4).ColumnWidth = 9.333333333333334
$ws.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws.Columns.Item(10).ColumnWidth = 8.0
$ws.Columns.Item(11).ColumnWidth = 17.666666666666668
$ws.Columns.Item(12).ColumnWidth = 18.5

# --- Row 21 (section title bar) -----------------------------------------
# Extend the title-bar formatting (style carried by J21) across the two
# new columns so the merged banner fully covers C21:L21.
$ws.Range("J21").Copy()
$ws.Range("K21:L21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Row 22 (column headers - internal names) ---------------------------
$ws.Range("J22").Copy()
$ws.Range("K22:L22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K22").Value = "classABeanFactory"
$ws.Range("L22").Value = "classBBeanFactory"

# --- Row 23 (column headers - display names) -----------------------------
$ws.Range("J23").Copy()
$ws.Range("K23:L23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K23").Value = "Class A Bean Factory"
$ws.Range("L23").Value = "Class B Bean Factory"

# --- Re-merge the section title bar across the new columns ---------------
$ws.Range("C21:L21").Merge()

# --- Row 25 (new data row, formatted like row 24) -------------------------
$ws.Range("C24:J24").Copy()
$ws.Range("C25:J25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Rows.Item(25).RowHeight = 15.75

$ws.Range("C25").Value = "A"
$ws.Range("D25").Value = "C"
$ws.Range("E25").Value = $false
$ws.Range("F25").Value = $true
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = $false
$ws.Range("I25").Value = "MM-dd-yyyy"
$ws.Range("J25").Value = $false

# --- Selection, matching the author's saved cursor position --------------
[void]$ws.Range("B24").Select()
